$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds numeric-looking text such as "216.71" or
# "0.510" that must stay plain text (matching the source feeds inline
# strings). Briefly forcing a Text number format before the write keeps
# Excel from auto-converting it to a Number (and dropping trailing
# zeroes); resetting the style back to Normal afterwards avoids leaving
# a stray cell format behind.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.828.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.28%  "

# Row 4
$ws.Range("E4").Value = "  -0.47%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.76%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.510"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.23%  "

# Row 7
$ws.Range("E7").Value = "  -0.45%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.254"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.38%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.01%  "

# Row 11
$ws.Range("E11").Value = "  -0.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.868.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.14%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.630.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.00%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.529"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.63%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.824.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.61%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.42%  "

# Row 20
$ws.Range("E20").Value = "  -0.58%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.69%  "

# Row 22
$ws.Range("E22").Value = "  +0.40%  "

# Row 23
$ws.Range("E23").Value = "  +2.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.41%  "

# Row 26
$ws.Range("E26").Value = "  -0.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.90%  "

# Row 28
$ws.Range("E28").Value = "  +0.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0501"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.32%  "

# Row 31
$ws.Range("E31").Value = "  -1.16%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.62%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.56"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.35%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.265.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.17%  "

# Row 36
$ws.Range("E36").Value = "  -0.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0177"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.69%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.533"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.08%  "

# Row 40
$ws.Range("E40").Value = "  -0.43%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.807"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.57%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.780.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "

# Row 44
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.20%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.10%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.17%  "

# Row 47
$ws.Range("E47").Value = "  -1.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.48%  "

# Row 49
$ws.Range("E49").Value = "  -0.69%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.86%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.52%  "
